$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new staff rows (Lab Instructors) below the existing data.
# Row 22: VIJAYA LAKSHMI
$ws.Cells.Item(22,1).Value = "VIJAYA LAKSHMI"
$ws.Cells.Item(22,2).Value = "Lab Instructor"
$ws.Cells.Item(22,3).Value = "/static/images/profile_photos/011/VEC-011-05-016.webp"
$ws.Cells.Item(22,10).Value = "VEC-011-05-016"

# Row 23: SURESH V
$ws.Cells.Item(23,1).Value = "SURESH V"
$ws.Cells.Item(23,2).Value = "Lab Instructor"
$ws.Cells.Item(23,3).Value = "/static/images/profile_photos/011/VEC-011-05-017.webp"
$ws.Cells.Item(23,10).Value = "VEC-011-05-017"

# Row 24: KANAGARAJ B
$ws.Cells.Item(24,1).Value = "KANAGARAJ B"
$ws.Cells.Item(24,2).Value = "Lab Instructor"
$ws.Cells.Item(24,3).Value = "/static/images/profile_photos/011/VEC-011-05-018.webp"
$ws.Cells.Item(24,10).Value = "VEC-011-05-018"

# Update the view state: scroll so row 11 is at the top and select C23,
# matching where the author ended up after inserting the new rows.
$win = $wb.Windows.Item(1)
$win.ScrollRow = 11
$win.ScrollColumn = 1
$ws.Range("C23").Select()

$wb.Save()
